$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2
    if ($text -eq $null) { continue }
    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }
    if ($parts[0].Equals("System")) { continue }
    if (($parts[0] -ieq "system") -or ($parts[1] -ieq "system")) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
